{"js": "// Append \" (Changed main)\" to the end of the first paragraph\n// (\"This is a Microsoft word document.\") as three separate runs:\n//   \" (\", \"Changed main\", \")\"\n// matching the target OOXML diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Insert each chunk right after the existing text, one at a time, so the\n// edits land in the same left-to-right order shown in the diff.\nfirstParagraph.insertText(\" (\", Word.InsertLocation.end);\nawait context.sync();\n\nfirstParagraph.insertText(\"Changed main\", Word.InsertLocation.end);\nawait context.sync();\n\nfirstParagraph.insertText(\")\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Append \" (Changed main)\" to the end of the first paragraph\n# (\"This is a Microsoft word document.\") as three separate\n# InsertAfter calls -- \" (\", \"Changed main\", \")\" -- matching the\n# target OOXML diff.\n$d = $word.ActiveDocument\n$r = $d.Paragraphs(1).Range\n\n# Paragraph ranges include the trailing paragraph mark, so trim it off\n# before appending to avoid spilling the new text into paragraph 2.\n$r.MoveEnd(1, -1) | Out-Null\n\n$r.Collapse(0)\n$r.InsertAfter(\" (\")\n\n$r.Collapse(0)\n$r.InsertAfter(\"Changed main\")\n\n$r.Collapse(0)\n$r.InsertAfter(\")\")\n"}
